$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 98.04003601545995
$ws.Range("B3").Value = 101.9599639845401
$ws.Range("B4").Value = 98.04003601545995
$ws.Range("B5").Value = 101.9599639845401
$ws.Range("B6").Value = 98.04003601545995
$ws.Range("B7").Value = 101.9599639845401
$ws.Range("B8").Value = 98.04003601545995
$ws.Range("B9").Value = 101.9599639845401
